# Auto-generated edit script applying numeric corrections to Kraken_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 681.28
$ws.Range("I6").Value = 387.07144
$ws.Range("K6").Value = 1161.21432
$ws.Range("M6").Value = -1049.21432

$ws.Range("H39").Value = 336
$ws.Range("I39").Value = 15.166667
$ws.Range("J39").Value = 1619.3334
$ws.Range("K39").Value = 45.500001
$ws.Range("L39").Value = 4858.0002
$ws.Range("M39").Value = 250.499999
$ws.Range("N39").Value = -5450.0002

$ws.Range("H132").Value = 2029.3478
$ws.Range("I132").Value = 2029.3478
$ws.Range("K132").Value = 6088.0434
$ws.Range("M132").Value = -3558.0434

$ws.Range("H137").Value = 1891.8334
$ws.Range("I137").Value = 1850.25
$ws.Range("J137").Value = 1975
$ws.Range("K137").Value = 5550.75
$ws.Range("L137").Value = 5925
$ws.Range("M137").Value = -3000.75
$ws.Range("N137").Value = -11025

$ws.Range("H138").Value = 3642.25
$ws.Range("I138").Value = 2036.625
$ws.Range("K138").Value = 6109.875
$ws.Range("M138").Value = -969.875

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2749.5
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 3499
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 3499
$ws.Range("M61").Value = -1788
$ws.Range("N61").Value = -3923

$ws.Range("H109").Value = 30377
$ws.Range("J109").Value = 30377
$ws.Range("L109").Value = 30377
$ws.Range("N109").Value = -33151

$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 3500
$ws.Range("K122").Value = 10500
$ws.Range("M122").Value = -8050

$ws.Range("H131").Value = 94997
$ws.Range("J131").Value = 94997
$ws.Range("L131").Value = 94997
$ws.Range("N131").Value = -105077

$ws.Range("H136").Value = 2749.5
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 3499
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 10497
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -15597

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 287
$ws.Range("I80").Value = 299.2857
$ws.Range("J80").Value = 265.5
$ws.Range("K80").Value = 299.2857
$ws.Range("L80").Value = 265.5
$ws.Range("M80").Value = 698.7143
$ws.Range("N80").Value = -2261.5

$ws.Range("H83").Value = 287
$ws.Range("I83").Value = 299.2857
$ws.Range("J83").Value = 265.5
$ws.Range("K83").Value = 1496.4285
$ws.Range("L83").Value = 1327.5
$ws.Range("M83").Value = 3495.5715
$ws.Range("N83").Value = -11311.5

$ws.Range("H32").Value = 1935.4445
$ws.Range("I32").Value = 1059.8572
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 1059.8572
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -743.8571999999999
$ws.Range("N32").Value = -5632

$ws.Range("H45").Value = 367
$ws.Range("I45").Value = 367
$ws.Range("K45").Value = 367
$ws.Range("M45").Value = 226

$ws.Range("H52").Value = 99995
$ws.Range("J52").Value = 99995
$ws.Range("L52").Value = 99995
$ws.Range("N52").Value = -100583

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 985.75
$ws.Range("J58").Value = 1249.6666
$ws.Range("L58").Value = 1249.6666
$ws.Range("N58").Value = -1655.6666

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0

$ws.Range("H136").Value = 985.75
$ws.Range("J136").Value = 1249.6666
$ws.Range("L136").Value = 3748.9998
$ws.Range("N136").Value = -8848.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1117.8889
$ws.Range("I7").Value = 1434.1428
$ws.Range("J7").Value = 11
$ws.Range("K7").Value = 4302.428400000001
$ws.Range("L7").Value = 33
$ws.Range("M7").Value = -4190.428400000001
$ws.Range("N7").Value = -257

$ws.Range("H112").Value = 1998.75
$ws.Range("I112").Value = 1331.6666
$ws.Range("J112").Value = 4000
$ws.Range("K112").Value = 3994.9998
$ws.Range("L112").Value = 12000
$ws.Range("M112").Value = -2886.9998
$ws.Range("N112").Value = -14216

$ws.Range("H132").Value = 1651.25
$ws.Range("I132").Value = 500
$ws.Range("J132").Value = 2035
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 18315
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -23375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4429.4287
$ws.Range("I122").Value = 3399.6
$ws.Range("K122").Value = 10198.8
$ws.Range("M122").Value = -7748.799999999999

$ws.Range("H132").Value = 6822.1113
$ws.Range("I132").Value = 6771.4287
$ws.Range("J132").Value = 6999.5
$ws.Range("K132").Value = 20314.2861
$ws.Range("L132").Value = 20998.5
$ws.Range("M132").Value = -17784.2861
$ws.Range("N132").Value = -26058.5

$ws.Range("H137").Value = 99995
$ws.Range("J137").Value = 99995
$ws.Range("L137").Value = 99995
$ws.Range("N137").Value = -110195

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3015.35
$ws.Range("J22").Value = 3220.4666
$ws.Range("L22").Value = 3220.4666
$ws.Range("N22").Value = -3810.4666

$ws.Range("H27").Value = 3015.35
$ws.Range("J27").Value = 3220.4666
$ws.Range("L27").Value = 3220.4666
$ws.Range("N27").Value = -3434.4666

$ws.Range("H32").Value = 13000
$ws.Range("I32").Value = 13000
$ws.Range("K32").Value = 13000
$ws.Range("M32").Value = -12683

$ws.Range("H40").Value = 12910.556
$ws.Range("I40").Value = 12910.556
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 12910.556
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -12774.556

$ws.Range("H95").Value = 18026
$ws.Range("J95").Value = 18026
$ws.Range("L95").Value = 18026
$ws.Range("N95").Value = -23518

$ws.Range("H132").Value = 4312.5

$ws.Range("H133").Value = 62500
$ws.Range("J133").Value = 62500
$ws.Range("L133").Value = 62500
$ws.Range("N133").Value = -67560

$ws.Range("H136").Value = 5173.25
$ws.Range("I136").Value = 4626.5713
$ws.Range("K136").Value = 13879.7139
$ws.Range("M136").Value = -11329.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4964
$ws.Range("I122").Value = 4624.6665
$ws.Range("K122").Value = 13873.9995
$ws.Range("M122").Value = -11423.9995

$ws.Range("H136").Value = 5714.125
$ws.Range("I136").Value = 4562.5
$ws.Range("K136").Value = 13687.5
$ws.Range("M136").Value = -11137.5
